$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply all updated cell values from the cryptos data refresh.
# Column D (Price) cells that look like plain decimal numbers need to be
# forced to Text so Excel does not auto-convert them to a Number and
# silently drop a significant trailing zero (e.g. "148.70" -> 148.7).
# NumberFormat is reset with ClearFormats() right after so the cell is
# left with no explicit style, matching the original inline-string cells.

$ws.Range('D2').Value = '72.327.71'
$ws.Range('E2').Value = '  +5.23%  '
$ws.Range('D3').Value = '4.079.80'
$ws.Range('E3').Value = '  +5.67%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '522.47'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.70'
$ws.Range('D6').ClearFormats()
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.726'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +19.80%  '
$ws.Range('D8').Value = '4.070.81'
$ws.Range('E8').Value = '  +5.65%  '
$ws.Range('E9').Value = '  +0.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.778'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +9.38%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.180'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +6.70%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000335'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +2.81%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '48.86'
$ws.Range('D13').ClearFormats()
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '11.07'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +8.83%  '
$ws.Range('D15').Value = '4.723.37'
$ws.Range('E15').Value = '  +5.80%  '
$ws.Range('D16').Value = '4.079.66'
$ws.Range('E16').Value = '  +5.50%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.57'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +4.39%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '21.38'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.25'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.08%  '
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').Value = '72.383.87'
$ws.Range('E21').Value = '  +5.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '448.60'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +6.82%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '104.05'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +19.63%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.63'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +6.82%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '15.05'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +7.36%  '
$ws.Range('E26').Value = '  +2.31%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.48'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.60%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.16'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +5.36%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '38.16'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +6.14%  '
$ws.Range('E30').Value = '  +2.91%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.31'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +16.92%  '
$ws.Range('E32').Value = '  +4.85%  '
$ws.Range('E33').Value = '  +5.58%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '683.44'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.17%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '67.89'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.12%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.67'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +13.55%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '42.52'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +6.98%  '
$ws.Range('D38').Value = '0.0₃0894'
$ws.Range('E38').Value = '  +5.35%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.435'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.19%  '
$ws.Range('E40').Value = '  +5.38%  '
$ws.Range('E41').Value = '  +9.20%  '
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0505'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +5.75%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.998'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.17%  '
$ws.Range('E45').Value = '  +1.61%  '
$ws.Range('E46').Value = '  +13.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.94'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +17.31%  '
$ws.Range('E48').Value = '  -0.50%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.41'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.07%  '
$ws.Range('B50').Value = 'FLOKI'
$ws.Range('C50').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000286'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +6.91%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.09'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +5.17%  '
